$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Replace the entire contents of a table cell with new text, preserving the
# cell's existing run formatting (Range.Text on the cell range keeps the
# surrounding run properties rather than inserting a brand new default run).
function Set-CellText($row, $col, $new) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $new
}

# Replace only the sub-range of a cell that starts at the first colon ":" in
# the cell (used for the "Crash Course: ..." cells, where the italic "Crash
# Course" run must stay untouched and only the trailing run changes).
function Set-CellSuffixAfterColon($row, $col, $new) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $fullText = $rng.Text
    $idx = $fullText.IndexOf(":")
    $subStart = $rng.Start + $idx
    $subEnd = $rng.End - 1
    $subRng = $d.Range($subStart, $subEnd)
    $subRng.Text = $new
}

# --- Unit III block (row 11-13) swaps with Unit IV block (row 15-17) ---

# Unit III header: "AI and the World" -> "Building AI: Business and Economics"
Set-CellText 11 3 "Building AI: Business and Economics"

# Week 6 crash course blurb
Set-CellSuffixAfterColon 12 3 ": Building AI: Business & Economics"

# Week 6 Thursday topic
Set-CellText 12 4 "Energy"

# Week 7 Tuesday topic
Set-CellText 13 3 "Labor Replacement I"

# Unit IV header: "Building AI: Business and Economics" -> "AI and the World"
Set-CellText 15 3 "AI and the World"

# Week 8 crash course blurb
Set-CellSuffixAfterColon 16 4 ": AI & The World"

# Week 9 Tuesday topic
Set-CellText 17 3 "AI Geopolitics"

# Week 11 Tuesday topic
Set-CellText 21 3 "Democracy & AI"

# Week 12 Thursday topic
Set-CellText 24 4 "Authoritarianism and AI"
